# Week 1 Assignment - Dungeon Crawler / Data Sheet.xlsx
# "Controls added. Improved UI"
#
# - Renames the original (only) sheet "Sheet1" -> "Rooms"
# - Adds a new "Items" sheet, positioned after "Rooms", and makes it active
# - Populates the Items sheet with an items table (Name / Description / LocationID)
#   bound to a Table (ListObject) mirroring the Room table layout
# - Clears the leftover "ID" number-format styling on Rooms!A2:A3 so the ID
#   column renders with the default style (matches the rest of the column)
# - Updates the selection on both sheets to match the new UI state

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Rename the existing sheet to "Rooms"
# ---------------------------------------------------------------------------
$rooms = $wb.Worksheets.Item(1)
$rooms.Name = "Rooms"

# ---------------------------------------------------------------------------
# 2. Add the new "Items" sheet right after "Rooms"
# ---------------------------------------------------------------------------
$items = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $rooms)
$items.Name = "Items"

# ---------------------------------------------------------------------------
# 3. Fill in the Items table headers + rows
# ---------------------------------------------------------------------------
$itemRows = @(
    @("Name", "Description", "LocationID"),
    @("Flashlight", "A black flashlight", 1),
    @("Crayon", "A red crayon", 2),
    @("Rope", "Some rope", 1),
    @("Weed", "A jar of weed", 4),
    @("Knife", "A karambit", 3),
    @("Banjo", "An old banjo with a missing string.", 3)
)

for ($r = 0; $r -lt $itemRows.Length; $r++) {
    $row = $itemRows[$r]
    $items.Cells.Item($r + 1, 1).Value = $row[0]
    $items.Cells.Item($r + 1, 2).Value = $row[1]
    $items.Cells.Item($r + 1, 3).Value = $row[2]
}

# Column widths matching the refreshed UI layout
$items.Columns.Item(1).ColumnWidth = 23.59
$items.Columns.Item(2).ColumnWidth = 31.92
$items.Columns.Item(3).ColumnWidth = 29.59

# First data row keeps the same "text" style used elsewhere in the workbook
# (mirrors how the Rooms sheet's imported rows are styled)
$items.Range("A2:B2").NumberFormat = "@"

# ---------------------------------------------------------------------------
# 4. Turn the item range into a table, like Rooms' Table1
# ---------------------------------------------------------------------------
$itemsTable = $items.ListObjects.Add(1, $items.Range("A1:C7"), $null, 1)
$itemsTable.Name = "Table2"

# ---------------------------------------------------------------------------
# 5. Clear the stray "ID" style on Rooms!A2:A3 (now rendered with default style)
# ---------------------------------------------------------------------------
$rooms.Range("A2").ClearFormats()
$rooms.Range("A3").ClearFormats()

# ---------------------------------------------------------------------------
# 6. Update selections / active sheet to match the saved UI state
# ---------------------------------------------------------------------------
[void]$rooms.Range("A2:K7").Select()
[void]$items.Range("C8").Select()
$items.Activate()
